# Changed the column names
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "Ghana"
$ws.Range("H1").Value = "Nigeria"
